# Purchases_Totals.xlsx update:
#  1. Move the "Peltier Cooler" line (bought by Yohan) from
#     "Juan's purchases" row 10 to "Yohan's purchases" row 2.
#  2. Fix up the totals formula on both of those sheets.
#  3. Add a new "Totals" sheet that rolls up each buyer's total and
#     computes an even split across the three purchasers.

$wb = $excel.ActiveWorkbook

$juan   = $wb.Worksheets.Item("Juan's purchases")
$robert = $wb.Worksheets.Item("Robert's purchases")
$yohan  = $wb.Worksheets.Item("Yohan's purchases")

# --- Move the Peltier Cooler row from Juan's sheet to Yohan's sheet ---
# Copy preserves both the values and the per-cell styles (e.g. the
# italic part-number cell, the currency-formatted price cell).
$juan.Range("A10:E10").Copy($yohan.Range("A2"))

# Remove the now-duplicated row from Juan's sheet; rows below shift up.
$juan.Rows.Item(10).Delete()

# Juan's total now only covers the remaining purchase rows (E2:E9)
# (selection on this sheet stays put at E12 from the original file)
$juan.Range("E11").Formula = "=SUM(E2:E9)"

# Yohan's total now sums the single purchase row just added
$yohan.Range("E10").Formula = "=SUM(E2)"
[void]$yohan.Range("E11").Select()

# --- Add the new "Totals" sheet, right after "Yohan's purchases" ---
# Adding it after $yohan both places it last and makes it the active
# sheet (matching activeTab on the workbook + tabSelected on the sheet).
$totals = $wb.Worksheets.Add($null, $yohan)
$totals.Name = "Totals"

$totals.Range("B3").Value = "Juan"
$totals.Range("C3").Formula = "='Juan''s purchases'!E11"

$totals.Range("B4").Value = "Robert"
$totals.Range("C4").Formula = "='Robert''s purchases'!E10"

$totals.Range("B5").Value = "Yohan"
$totals.Range("C5").Formula = "='Yohan''s purchases'!E10"

$totals.Range("B6").Value = "Total"
$totals.Range("C6").Formula = "=SUM(C3:C5)"

$totals.Range("B8").Value = "Even per person"
$totals.Range("C8").Formula = "=C6/3"

$totals.Columns.Item(2).ColumnWidth = 13

[void]$totals.Range("B10").Select()
